# Update the "想去人数" (interest count) column F values on the
# "展览" sheet and the "全部类型" sheet to reflect newly generated
# output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1) - column F rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 386
$wsExhibit.Range("F4").Value  = 443
$wsExhibit.Range("F5").Value  = 27
$wsExhibit.Range("F7").Value  = 252
$wsExhibit.Range("F8").Value  = 13924
$wsExhibit.Range("F9").Value  = 88
$wsExhibit.Range("F10").Value = 83
$wsExhibit.Range("F11").Value = 5618
$wsExhibit.Range("F12").Value = 578
$wsExhibit.Range("F13").Value = 48
$wsExhibit.Range("F14").Value = 36
$wsExhibit.Range("F16").Value = 1222
$wsExhibit.Range("F19").Value = 752
$wsExhibit.Range("F20").Value = 2905
$wsExhibit.Range("F22").Value = 10405
$wsExhibit.Range("F23").Value = 1187
$wsExhibit.Range("F25").Value = 41
$wsExhibit.Range("F26").Value = 3705

# Sheet "全部类型" (sheetId 4) - same events, offset by extra rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 386
$wsAll.Range("F5").Value  = 443
$wsAll.Range("F6").Value  = 27
$wsAll.Range("F8").Value  = 252
$wsAll.Range("F9").Value  = 13924
$wsAll.Range("F10").Value = 88
$wsAll.Range("F11").Value = 83
$wsAll.Range("F12").Value = 5618
$wsAll.Range("F13").Value = 578
$wsAll.Range("F14").Value = 48
$wsAll.Range("F15").Value = 36
$wsAll.Range("F17").Value = 1222
$wsAll.Range("F20").Value = 752
$wsAll.Range("F21").Value = 2905
$wsAll.Range("F24").Value = 10405
$wsAll.Range("F25").Value = 1187
$wsAll.Range("F27").Value = 41
$wsAll.Range("F28").Value = 3705
